$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2, 2, 2, 1, 1, 1, 2, 1, 2, 2, 2, 1, 2, 2)

for ($i = 0; $i -lt $values.Length; $i++) {
    # Column B is index 2, so offset by 2
    $ws.Cells.Item(2, $i + 2).Value = $values[$i]
}
